$d = $word.ActiveDocument

# Locate the existing "Erstellen der GitRepository" bullet; the new bullets
# are appended right after it.
$lastPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Erstellen der GitRepository*") {
        $lastPara = $p
    }
}

$newItems = @(
    "Informationsbeschaffung zu neuen Klassen",
    "Erstellen des JFrame",
    "Erstellen des JPanel",
    "Erstellen der Zeichenfläche mit paintComponent und Konstruktor Panel"
)

foreach ($text in $newItems) {
    # InsertParagraphAfter() clones the current paragraph's formatting
    # (style "Listenabsatz" + the w:numPr list numbering) onto the new one.
    $lastPara.Range.InsertParagraphAfter()
    $count = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs($count)
    $lastPara.Range.Text = $text
}
